# Admin Module - Debenture_Type_setting first push
#
# 1) Update 4 data values on "AccountOpening_Loan_ACOPL_TwoWh" (sheet12)
# 2) Adjust its sheetView (scroll/selection)
# 3) Clone that sheet to create the new "Debenture_Type_setting" sheet,
#    re-point its header/data cells, and tidy up formatting/selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Existing sheet: AccountOpening_Loan_ACOPL_TwoWh
# ---------------------------------------------------------------------
$loanSheet = $wb.Worksheets.Item("AccountOpening_Loan_ACOPL_TwoWh")

$loanSheet.Range("E2").Value = 102000000001
$loanSheet.Range("F2").Value = "A"
$loanSheet.Range("AE2").Value = "25/12/2300"
$loanSheet.Range("AF2").Value = "M"

# ---------------------------------------------------------------------
# 2) Create the new sheet by copying the loan sheet (keeps styles,
#    column widths, row heights, number formats, etc.) then rename it.
# ---------------------------------------------------------------------
$loanSheet.Copy($null, $loanSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Debenture_Type_setting"

# Row 1 headers that differ from the source sheet.
$newSheet.Range("E1").Value = "Name"
$newSheet.Range("F1").Value = "validAmount"
$newSheet.Range("G1").Value = "Duration"
$newSheet.Range("H1").Value = "minimum"
$newSheet.Range("I1").Value = "maximum"

# Blank out the remaining header cells (style stays, text goes).
$newSheet.Range("J1:AH1").ClearContents()

# Row 2 sample data.
$newSheet.Range("A2").Value = "Debenture_Type_setting"
$newSheet.Range("E2").Value = "NewDebe"

# F2/H2 need to look like plain numbers (copy number format off the
# neighbouring already-general-formatted cells) before assigning values.
$newSheet.Range("G2").Copy()
$newSheet.Range("F2").PasteSpecial(-4122)
$newSheet.Range("I2").Copy()
$newSheet.Range("H2").PasteSpecial(-4122)

$newSheet.Range("F2").Value = 1500
$newSheet.Range("G2").Value = 12
$newSheet.Range("H2").Value = 10
$newSheet.Range("I2").Value = 100

# Clear the rest of row 2 (style stays where it existed, AG2/AH2 had no
# style to begin with so they disappear entirely, matching the target).
$newSheet.Range("J2:AH2").ClearContents()

# Column A is a touch wider on the new sheet than on the source sheet.
$newSheet.Columns("A").ColumnWidth = 13.67

# Row 2 is shorter on the new sheet (less text wrapping needed).
$newSheet.Rows(2).RowHeight = 29

# New sheet's view: not the tab-selected sheet, simple top-left selection.
$newSheet.Range("G9").Select()

# ---------------------------------------------------------------------
# 3) Restore focus/scroll-selection on the original sheet.
# ---------------------------------------------------------------------
$loanSheet.Activate()
$loanSheet.Range("AD7").Select()
